# Append the 2026/01/31 data row to the "ModCounts" sheet.
#
# The new row must land on row 82 with the same cell style (s="1",
# center/center alignment) used by every other data row, and column A must
# store the date as a literal text string "2026/01/31" (the sheet has no
# date number-format anywhere — every existing "Date" cell is plain text),
# not get auto-converted into an Excel serial date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 81
$newRow = $lastRow + 1

$srcRange = $ws.Range("A$lastRow" + ":C$lastRow")
$dstRange = $ws.Range("A$newRow" + ":C$newRow")

# Clone the previous row first: this carries the existing style (s="1",
# centered alignment) onto every cell of the new row without us having to
# touch styles.xml at all.
$srcRange.Copy($dstRange)

# Stage the date text on a scratch cell, far away from the used range, as a
# formula result. A formula-computed string is never re-interpreted as a
# date by Excel, so it stays plain text with General format - no new style
# gets created. We then paste just that *value* onto A82 so its existing
# style (already copied above) is left untouched.
$scratch = $ws.Cells.Item(1, 26)
$scratch.Formula = '="2026/01/31"'
$scratch.Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)  # xlPasteValues
$scratch.Formula = ""

$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1168
